$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-WrapForEach($addrCsv) {
    foreach ($a in $addrCsv.Split(',')) {
        $ws.Range($a).WrapText = $true
    }
}

# 1) "Numele și Prenumele:" label gains a two-space prefix (A3).
$ws.Range("A3").Value = "  Numele și Prenumele:"

# 2) Column B gets narrower (250.71.. -> 200.71.. characters).
$ws.Range("B:B").ColumnWidth = 200.7109375

# 3) New per-sentence word-count values, column A.
$counts = @{7=10; 19=10; 31=8; 43=10; 55=9; 67=9; 79=6; 91=10; 103=10; 115=2; 127=8; 139=10}
foreach ($row in $counts.Keys) {
    $ws.Cells.Item($row, 1).Value = $counts[$row]
}

# 4) The "filler word count" cells in column A (style: fontId 1, no fill):
#    - font colour goes from gray FF808080 to near-white FFE3E3E3
#    - alignment becomes centered both ways
#    This covers the existing label/number cells plus the freshly added
#    per-sentence word-count cells above — every address is touched
#    individually so the COM layer doesn't silently drop all but the
#    first area of a multi-area Range.
$labelRows = @(1,6,18,30,42,54,66,78,90,102,114,126,138,7,19,31,43,55,67,79,91,103,115,127,139)
foreach ($row in $labelRows) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Font.Color = 14935011
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# 5) Word-wrap turned on for the colored section rows:
#    - A5/A17/... (red "word" header cells) also gain horizontal centering
$redHeaderRows = @(5,17,29,41,53,65,77,89,101,113,125,137)
foreach ($row in $redHeaderRows) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.HorizontalAlignment = -4108
    $cell.WrapText = $true
}

#    - B5/B17/... (orange "definition" header cells)
foreach ($row in $redHeaderRows) {
    $ws.Cells.Item($row, 2).WrapText = $true
}

#    - B6,B8,... (empty blue filler-sentence cells)
$blue = "B6,B8,B10,B12,B14,B18,B20,B22,B24,B26,B32,B34,B36,B38,B42,B44,B46,B48,B50," +
        "B56,B58,B60,B62,B68,B70,B72,B74,B82,B84,B86,B90,B92,B94,B96,B98,B102,B104," +
        "B106,B108,B110,B122,B128,B130,B132,B134,B138,B140,B142,B144,B146"
Set-WrapForEach $blue

#    - B7,B9,... (empty light-blue filler-sentence cells)
$lightBlue = "B7,B9,B11,B13,B15,B19,B21,B23,B25,B27,B33,B35,B37,B39,B43,B45,B47,B49,B51," +
             "B55,B57,B59,B61,B63,B67,B69,B71,B73,B75,B83,B85,B87,B91,B93,B95,B97,B99," +
             "B103,B105,B107,B109,B111,B123,B129,B131,B133,B135,B139,B141,B143,B145,B147"
Set-WrapForEach $lightBlue

#    - B30,B54,... (filled green example-sentence cells)
Set-WrapForEach "B30,B54,B66,B78,B80,B114,B116,B118,B120,B126"

#    - B31,B79,... (filled light-green example-sentence cells)
Set-WrapForEach "B31,B79,B81,B115,B117,B119,B121,B127"
